# Applies the "Updated cryptos list" refresh: new Price (D) and
# Volume(1h) (E) quotes for each coin row, plus a swap of the
# PancakeSwap / Decentraland rows (46-47, incl. B and C columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = '27.528.07'
$ws.Cells.Item(2, 5).Value = '  -1.94%  '

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = '1.750.78'
$ws.Cells.Item(3, 5).Value = '  -2.18%  '

# Row 4: TetherUSD
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.001'
$ws.Cells.Item(4, 5).Value = '  -0.06%  '

# Row 5: BNB
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '324.20'
$ws.Cells.Item(5, 5).Value = '  +0.20%  '

# Row 6: USDC
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.9999'
$ws.Cells.Item(6, 5).Value = '  -0.04%  '

# Row 7: XRP
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4465'
$ws.Cells.Item(7, 5).Value = '  +3.53%  '

# Row 8: Cardano
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3605'
$ws.Cells.Item(8, 5).Value = '  -0.50%  '

# Row 9: Dogecoin
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07512'
$ws.Cells.Item(9, 5).Value = '  +0.04%  '

# Row 10: OKB
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '41.93'
$ws.Cells.Item(10, 5).Value = '  -6.22%  '

# Row 11: Polygon
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1.093'
$ws.Cells.Item(11, 5).Value = '  -1.90%  '

# Row 12: BinanceUSD
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.9997'
$ws.Cells.Item(12, 5).Value = '  -0.07%  '

# Row 13: Solana
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '20.63'
$ws.Cells.Item(13, 5).Value = '  -4.56%  '

# Row 14: Polkadot
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.033'
$ws.Cells.Item(14, 5).Value = '  -1.86%  '

# Row 15: Chainlink
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.129'
$ws.Cells.Item(15, 5).Value = '  -2.55%  '

# Row 16: WrappedEther
$ws.Cells.Item(16, 4).Value = '1.751.10'
$ws.Cells.Item(16, 5).Value = '  -2.08%  '

# Row 17: Litecoin
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '93.38'
$ws.Cells.Item(17, 5).Value = '  +1.27%  '

# Row 18: ShibaInu
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.00001062'
$ws.Cells.Item(18, 5).Value = '  -0.31%  '

# Row 19: TRON
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06400'
$ws.Cells.Item(19, 5).Value = '  +0.68%  '

# Row 20: Dai
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '1.000'
$ws.Cells.Item(20, 5).Value = '  +0.01%  '

# Row 21: Avalanche
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '16.80'
$ws.Cells.Item(21, 5).Value = '  -2.51%  '

# Row 22: Uniswap
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.853'
$ws.Cells.Item(22, 5).Value = '  -1.98%  '

# Row 23: WrappedBTC
$ws.Cells.Item(23, 4).Value = '27.568.62'
$ws.Cells.Item(23, 5).Value = '  -1.80%  '

# Row 24: Cosmos
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '11.17'
$ws.Cells.Item(24, 5).Value = '  -1.98%  '

# Row 25: Toncoin
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.105'
$ws.Cells.Item(25, 5).Value = '  -1.70%  '

# Row 26: Monero
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '161.62'
$ws.Cells.Item(26, 5).Value = '  +1.66%  '

# Row 27: EthereumClassic
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '20.48'
$ws.Cells.Item(27, 5).Value = '  +0.61%  '

# Row 28: WrappedliquidstakedEther2.0
$ws.Cells.Item(28, 4).Value = '1.951.17'

# Row 29: LidoDAOToken
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.093'
$ws.Cells.Item(29, 5).Value = '  -3.53%  '

# Row 30: BitcoinCash
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '125.26'
$ws.Cells.Item(30, 5).Value = '  -1.54%  '

# Row 31: ImmutableX
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.084'
$ws.Cells.Item(31, 5).Value = '  -6.18%  '

# Row 32: HuobiToken
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.660'
$ws.Cells.Item(32, 5).Value = '  +3.92%  '

# Row 33: Stellar
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.09014'
$ws.Cells.Item(33, 5).Value = '  +0.40%  '

# Row 34: Filecoin
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.553'
$ws.Cells.Item(34, 5).Value = '  -3.45%  '

# Row 35: Aptos
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '11.99'
$ws.Cells.Item(35, 5).Value = '  -4.87%  '

# Row 36: VeChain
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.02298'
$ws.Cells.Item(36, 5).Value = '  -0.76%  '

# Row 37: Hedera
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.06021'
$ws.Cells.Item(37, 5).Value = '  -0.52%  '

# Row 38: Algorand
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.2086'
$ws.Cells.Item(38, 5).Value = '  -1.26%  '

# Row 39: TheSandbox
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.6346'
$ws.Cells.Item(39, 5).Value = '  -1.50%  '

# Row 40: InternetComputer(DFINITY)
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '4.957'
$ws.Cells.Item(40, 5).Value = '  -2.62%  '

# Row 41: TrustWalletToken
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.206'
$ws.Cells.Item(41, 5).Value = '  +1.51%  '

# Row 42: Frax
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.9997'
$ws.Cells.Item(42, 5).Value = '  +0.01%  '

# Row 43: WEMIXTOKEN
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.380'
$ws.Cells.Item(43, 5).Value = '  -2.97%  '

# Row 44: FraxShare
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '7.779'
$ws.Cells.Item(44, 5).Value = '  -0.72%  '

# Row 45: EnergySwap
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '13.18'
$ws.Cells.Item(45, 5).Value = '  -2.79%  '

# Row 46: PancakeSwap
$ws.Cells.Item(46, 2).Value = 'Decentraland'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.5888'
$ws.Cells.Item(46, 5).Value = '  -1.68%  '

# Row 47: Decentraland
$ws.Cells.Item(47, 2).Value = 'PancakeSwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.714'
$ws.Cells.Item(47, 5).Value = '  +0.21%  '

# Row 48: Quant
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '121.96'
$ws.Cells.Item(48, 5).Value = '  -2.16%  '

# Row 49: NEARProtocol
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.956'
$ws.Cells.Item(49, 5).Value = '  -1.19%  '

# Row 50: EOS
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.149'
$ws.Cells.Item(50, 5).Value = '  -0.33%  '

# Row 51: Cronos
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.06857'
$ws.Cells.Item(51, 5).Value = '  -1.29%  '
